$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (matches the workbook's original inline-string / text representation)
    # without altering the cell's number format.
    $ws.Range($addr).Value = "'" + $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "97.523.83"
Set-TextValue "E2" "  +2.77%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.599.11"
Set-TextValue "E3" "  +1.03%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.08%  "

# Row 5 - Solana
Set-TextValue "D5" "242.52"
Set-TextValue "E5" "  +2.72%  "

# Row 6 - XRP
Set-TextValue "E6" "  +17.59%  "

# Row 7 - BNB
Set-TextValue "D7" "653.78"
Set-TextValue "E7" "  -0.11%  "

# Row 8 - Dogecoin
Set-TextValue "E8" "  +9.68%  "

# Row 9 - USDC
Set-TextValue "E9" "  -0.02%  "

# Row 10 - Cardano
Set-TextValue "E10" "  +4.81%  "

# Row 11 - LidoStakedEther
Set-TextValue "D11" "3.597.52"
Set-TextValue "E11" "  +0.96%  "

# Row 12 - Avalanche
Set-TextValue "D12" "44.51"
Set-TextValue "E12" "  +5.48%  "

# Row 13 - TRON
Set-TextValue "E13" "  +1.02%  "

# Row 14 - Toncoin
Set-TextValue "E14" "  +0.65%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.267.99"
Set-TextValue "E15" "  +1.08%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "97.285.27"

# Row 17 - ShibaInu
Set-TextValue "E17" "  +3.84%  "

# Row 18 - Polkadot
Set-TextValue "D18" "8.68"
Set-TextValue "E18" "  +2.14%  "

# Row 19 - WrappedEther
Set-TextValue "D19" "3.599.76"
Set-TextValue "E19" "  +1.18%  "

# Row 20 - Uniswap
Set-TextValue "D20" "12.57"
Set-TextValue "E20" "  -0.85%  "

# Row 21 - Chainlink
Set-TextValue "D21" "18.14"
Set-TextValue "E21" "  +2.40%  "

# Row 22 - Stellar
Set-TextValue "D22" "0.530"
Set-TextValue "E22" "  +10.33%  "

# Row 23 - SuiNetwork
Set-TextValue "D23" "3.49"
Set-TextValue "E23" "  +1.42%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "518.47"
Set-TextValue "E24" "  +2.21%  "

# Row 25 - PEPE
Set-TextValue "D25" "0.0000212"
Set-TextValue "E25" "  +8.69%  "

# Row 26 - NEARProtocol
Set-TextValue "E26" "  +2.44%  "

# Row 27 - Litecoin
Set-TextValue "D27" "102.41"
Set-TextValue "E27" "  +8.03%  "

# Row 28 - Aptos
Set-TextValue "D28" "13.17"
Set-TextValue "E28" "  +4.96%  "

# Row 29 - WrappedeETH
Set-TextValue "D29" "3.792.44"
Set-TextValue "E29" "  +1.07%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.171"
Set-TextValue "E30" "  +19.01%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "12.08"
Set-TextValue "E31" "  +5.50%  "

# Row 32 - PancakeSwap
Set-TextValue "E32" "  -1.02%  "

# Row 33 - Dai
Set-TextValue "E33" "  +0.17%  "

# Row 34 - Cronos
Set-TextValue "E34" "  +6.11%  "

# Row 35 - Binance-PegBSC-USD
Set-TextValue "E35" "  +0.28%  "

# Row 36 - EthereumClassic
Set-TextValue "D36" "31.98"
Set-TextValue "E36" "  +0.80%  "

# Row 37 - Bittensor
Set-TextValue "D37" "619.05"
Set-TextValue "E37" "  +6.88%  "

# Row 38 - PolygonEcosystemToken
Set-TextValue "D38" "0.573"
Set-TextValue "E38" "  +3.56%  "

# Row 39 - RenderToken
Set-TextValue "D39" "8.74"
Set-TextValue "E39" "  +3.44%  "

# Row 40 - Fetch.AI
Set-TextValue "E40" "  -2.88%  "

# Row 41 - Kaspa
Set-TextValue "E41" "  +3.26%  "

# Row 42 - ImmutableX
Set-TextValue "E42" "  +7.85%  "

# Row 43 - ARBITRUM
Set-TextValue "E43" "  +3.16%  "

# Row 44 - USDe
Set-TextValue "E44" "  -0.01%  "

# Row 45 - Filecoin
Set-TextValue "D45" "6.03"
Set-TextValue "E45" "  +5.20%  "

# Row 46 and 47 - VeChain/Algorand swapped order
Set-TextValue "B46" "Algorand"
Set-TextValue "C46" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D46" "0.434"
Set-TextValue "E46" "  +43.21%  "

Set-TextValue "B47" "VeChain"
Set-TextValue "C47" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0442"
Set-TextValue "E47" "  +7.25%  "

# Row 48 - Stacks
Set-TextValue "E48" "  +2.04%  "

# Row 49 - WhiteBITCoin
Set-TextValue "D49" "23.65"
Set-TextValue "E49" "  +1.13%  "

# Row 50 - Cosmos
Set-TextValue "E50" "  +5.58%  "

# Row 51 - dogwifhat
Set-TextValue "D51" "3.31"
Set-TextValue "E51" "  +8.41%  "
